$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2976
$ws1.Range("F3").Value = 6420
$ws1.Range("F6").Value = 537
$ws1.Range("F7").Value = 67
$ws1.Range("F9").Value = 2935
$ws1.Range("F11").Value = 40
$ws1.Range("F12").Value = 7505
$ws1.Range("F19").Value = 9159
$ws1.Range("F30").Value = 112
$ws1.Range("F33").Value = 2617
$ws1.Range("F36").Value = 173
$ws1.Range("F38").Value = 768
$ws1.Range("F39").Value = 3920
$ws1.Range("F40").Value = 211
$ws1.Range("F41").Value = 39
$ws1.Range("F43").Value = 90
$ws1.Range("F44").Value = 27
$ws1.Range("F45").Value = 243

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 267
$ws2.Range("F7").Value = 128
$ws2.Range("F8").Value = 30
$ws2.Range("F18").Value = 168

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2976
$ws4.Range("F5").Value = 267
$ws4.Range("F6").Value = 6420
$ws4.Range("F8").Value = 128
$ws4.Range("F10").Value = 537
$ws4.Range("F11").Value = 67
$ws4.Range("F13").Value = 2935
$ws4.Range("F15").Value = 30
$ws4.Range("F17").Value = 40
$ws4.Range("F18").Value = 7505
$ws4.Range("F24").Value = 9160
$ws4.Range("F31").Value = 112
$ws4.Range("F34").Value = 2617
$ws4.Range("F36").Value = 173
$ws4.Range("F38").Value = 769
$ws4.Range("F39").Value = 168
$ws4.Range("F40").Value = 3920
$ws4.Range("F41").Value = 211
$ws4.Range("F42").Value = 39
$ws4.Range("F45").Value = 90
$ws4.Range("F46").Value = 243
